$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 121, shifting the existing rows
# 121-159 down to 123-161 (preserving all of their data/formatting).
$ws.Rows.Item(121).Resize(2).Insert()

# Row 121: new "Primera" quality record dated 44988
$ws.Range("A121").Value = 1
$ws.Range("B121").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C121").Value = "Arica y Parinacota"
$ws.Range("D121").Value = 44988
$ws.Range("E121").Value = 15
$ws.Range("F121").Value = 100112042
$ws.Range("G121").Value = "Locoto"
$ws.Range("H121").Value = "Sin especificar"
$ws.Range("I121").Value = "Primera"
$ws.Range("J121").Value = 100
$ws.Range("K121").Value = 48000
$ws.Range("L121").Value = 50000
$ws.Range("M121").Value = 49000
$ws.Range("N121").Value = "$/caja 20 kilos"
$ws.Range("O121").Value = "Región de Arica y Parinacota"
$ws.Range("P121").Value = 2450
$ws.Range("Q121").Value = 20
$ws.Range("R121").Value = "Hortaliza"

# Row 122: new "Segunda" quality record dated 44988
$ws.Range("A122").Value = 1
$ws.Range("B122").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C122").Value = "Arica y Parinacota"
$ws.Range("D122").Value = 44988
$ws.Range("E122").Value = 15
$ws.Range("F122").Value = 100112042
$ws.Range("G122").Value = "Locoto"
$ws.Range("H122").Value = "Sin especificar"
$ws.Range("I122").Value = "Segunda"
$ws.Range("J122").Value = 120
$ws.Range("K122").Value = 28000
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = 29000
$ws.Range("N122").Value = "$/caja 20 kilos"
$ws.Range("O122").Value = "Región de Arica y Parinacota"
$ws.Range("P122").Value = 1450
$ws.Range("Q122").Value = 20
$ws.Range("R122").Value = "Hortaliza"
